$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for new "Save" column - copy formatting (bold, border, centered) from G1
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("H1").Value = "Save"

# Save values per row (2-35)
$saveValues = @{
    2  = 0
    3  = 1
    4  = 0
    5  = 1
    6  = 0
    7  = 1
    8  = 0
    9  = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 0
    27 = 0
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
}

foreach ($row in $saveValues.Keys) {
    $ws.Range("H$row").Value = $saveValues[$row]
}
